$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated values: each card's info combined into a single
# Python tuple-repr string per row (A2:A10), replacing the old
# one-attribute-per-row layout (A2:A38).
$values = @(
    "('Bribery', ['{3}{U}{U}', 'Sorcery', 'Search target opponent" + [char]0x2019 + "s library for a creature card and put that card onto the battlefield under your control. Then that player shuffles their library.'])",
    "('Crucible of Worlds', ['{3}', 'Artifact', 'You may play lands from your graveyard.'])",
    "('Genesis', ['{4}{G}', 'Creature " + [char]0x2014 + " Incarnation', 'At the beginning of your upkeep, if Genesis is in your graveyard, you may pay {2}{G}. If you do, return target creature card from your graveyard to your hand.', '4/4'])",
    "('Golem', ['Token Enchantment Artifact Creature " + [char]0x2014 + " Golem', '3/3'])",
    "('Imperial Recruiter', ['{2}{R}', 'Creature " + [char]0x2014 + " Human Advisor', 'When Imperial Recruiter enters the battlefield, search your library for a creature card with power 2 or less, reveal it, put it into your hand, then shuffle your library.', '1/1'])",
    "('Overwhelming Forces', ['{6}{B}{B}', 'Sorcery', 'Destroy all creatures target opponent controls. Draw a card for each creature destroyed this way.'])",
    "('Show and Tell', ['{2}{U}', 'Sorcery', 'Each player may put an artifact, creature, enchantment, or land card from their hand onto the battlefield.'])",
    "('Swords to Plowshares', ['{W}', 'Instant', 'Exile target creature. Its controller gains life equal to its power.'])",
    "('Vindicate', ['{1}{W}{B}', 'Sorcery', 'Destroy target permanent.'])"
)

# Clear out the old rows beyond the new range (A11:A38) first.
$ws.Range("A11:A38").ClearContents()

# Write the new consolidated values into A2:A10.
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
